# "fully tested alu (otherthan add/sub)"
# Mark ALU ops tested in column G ("Tested and Working?") for rows 3-15 of
# Sheet2, and clear the now-redundant "Finished?" marks (column F) for the
# two rows (17/18 - slt/sub) whose testing status regressed/was corrected.

$wb = $excel.ActiveWorkbook

# ---- Sheet1: restore the scroll position of the frozen pane ----
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws1.Activate()
$win = $wb.Windows.Item(1)
$win.ScrollRow = 6
$win.ScrollColumn = 4

# ---- Sheet2: update the ALU test-tracking table ----
$ws2 = $wb.Worksheets.Item("Sheet2")
$ws2.Activate()

# Mark rows 3-15 (all ALU ops other than slt/sub handling in 17/18) as
# "Tested and Working?" = x
$ws2.Range("G3:G15").Value = "x"

# These two ops (row 17 = slt, row 18 = sub) are no longer marked
# "Finished?" in column F.
$ws2.Range("F17").ClearContents()
$ws2.Range("F18").ClearContents()

# Restore the active selection left on the sheet.
$ws2.Range("E19").Select()
